# Generate Report for Handoff
# Adds two new files (8b727eba-5984-4079-9ee2-9297b5953605.md and
# 9060964e-651e-4254-839c-9f639946a38a.md) as new rows to the Overview,
# zh-cn and de-de report tables.

$wb = $excel.ActiveWorkbook

$file1Name = "8b727eba-5984-4079-9ee2-9297b5953605.md"
$file1Path = "e2e\8b727eba-5984-4079-9ee2-9297b5953605.md"
$file2Name = "9060964e-651e-4254-839c-9f639946a38a.md"
$file2Path = "e2e\9060964e-651e-4254-839c-9f639946a38a.md"

$dateTime = "2016-08-20 04:44:39"
$zhDateTime = "2016-08-20 04:44:35"
$deDateTime = "2016-08-20 04:44:39"
$zeroDate = "0001-01-01 00:00:00"

$zhXlf1 = "8b727eba-5984-4079-9ee2-9297b5953605.8e23c20f8deb1436dc7f8251f45d8b2594291415.zh-cn.xlf"
$zhXlf2 = "9060964e-651e-4254-839c-9f639946a38a.cba59315bd9bfecd72cf719c92bbd38d66374010.zh-cn.xlf"
$deXlf1 = "8b727eba-5984-4079-9ee2-9297b5953605.8e23c20f8deb1436dc7f8251f45d8b2594291415.de-de.xlf"
$deXlf2 = "9060964e-651e-4254-839c-9f639946a38a.cba59315bd9bfecd72cf719c92bbd38d66374010.de-de.xlf"

# ============================================================
# Sheet "Overview"
# ============================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4 - 8b727eba...
$wsOverview.Range("A4").Value = $file1Name
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = $dateTime
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file1Name", "", "", $file1Path) | Out-Null

# Row 5 - 9060964e...
$wsOverview.Range("A5").Value = $file2Name
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = $dateTime
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file2Name", "", "", $file2Path) | Out-Null

# ============================================================
# Sheet "zh-cn"
# ============================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

# Row 4 - 8b727eba...
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = $zhXlf1
$wsZh.Range("H4").Value = $zhDateTime
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $zeroDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file1Name", "", "", $file1Name) | Out-Null

# Row 5 - 9060964e...
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = $zhXlf2
$wsZh.Range("H5").Value = $zhDateTime
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $zeroDate
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file2Name", "", "", $file2Name) | Out-Null

# ============================================================
# Sheet "de-de"
# ============================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

# Row 4 - 8b727eba...
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = $deXlf1
$wsDe.Range("H4").Value = $deDateTime
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $zeroDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file1Name", "", "", $file1Name) | Out-Null

# Row 5 - 9060964e...
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = $deXlf2
$wsDe.Range("H5").Value = $deDateTime
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $zeroDate
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/$file2Name", "", "", $file2Name) | Out-Null
